$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: region names reshuffled among rows 3, 4, 5, 7 -------------
$ws.Range("A3").Value = "São Paulo"
$ws.Range("A4").Value = "Goiás"
$ws.Range("A5").Value = "Rio Grande do Sul"
$ws.Range("A7").Value = "Distrito Federal"

# --- Column C: quarter date text, updated for every data row (2-10) -----
# Use a text number format while writing so Excel doesn't auto-convert the
# "dd/mm/yyyy"-looking string into a date serial, then drop back to the
# default "Normal" style so no residual formatting is left on the cells.
$cRng = $ws.Range("C2:C10")
$cRng.NumberFormat = "@"
$cRng.Value = "01/07/2024"
$cRng.Style = "Normal"

# --- Column D: updated numeric values for every data row (2-10) ---------
$ws.Range("D2").Value = 56.37035083091533
$ws.Range("D3").Value = 55.26409595101801
$ws.Range("D4").Value = 54.49565798263193
$ws.Range("D5").Value = 54.4818048232345
$ws.Range("D6").Value = 54.41783649876135
$ws.Range("D7").Value = 54.28482393268931
$ws.Range("D8").Value = 46.78631051752922
$ws.Range("D9").Value = 43.92015762871584
$ws.Range("D10").Value = 50.66514405698735
